$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure text-formatted columns (B:E) keep their literal text representation
# (e.g. "26.521.79", "0.9997", "  +0.03%  ") instead of being auto-converted
# to numbers by Excel when assigned.
$ws.Range("B2:E51").NumberFormat = "@"

$ws.Range('D2').Value = '26.521.79'
$ws.Range('E2').Value = '  +0.03%  '
$ws.Range('D3').Value = '1.736.77'
$ws.Range('E3').Value = '  +0.09%  '
$ws.Range('E4').Value = '  -0.07%  '
$ws.Range('D5').Value = '246.89'
$ws.Range('E5').Value = '  +1.23%  '
$ws.Range('D6').Value = '0.9997'
$ws.Range('D7').Value = '0.4895'
$ws.Range('E7').Value = '  +2.19%  '
$ws.Range('D8').Value = '0.2666'
$ws.Range('E8').Value = '  +0.11%  '
$ws.Range('D9').Value = '0.06367'
$ws.Range('E9').Value = '  +2.31%  '
$ws.Range('D10').Value = '1.731.00'
$ws.Range('E10').Value = '  -0.24%  '
$ws.Range('D11').Value = '0.07040'
$ws.Range('E11').Value = '  -1.25%  '
$ws.Range('D12').Value = '15.71'
$ws.Range('E12').Value = '  -0.03%  '
$ws.Range('D13').Value = '4.601'
$ws.Range('E13').Value = '  +1.63%  '
$ws.Range('D14').Value = '0.6101'
$ws.Range('E14').Value = '  -0.60%  '
$ws.Range('D15').Value = '77.43'
$ws.Range('E15').Value = '  +0.84%  '
$ws.Range('D16').Value = '0.9997'
$ws.Range('E16').Value = '  -0.10%  '
$ws.Range('D17').Value = '0.000007422'
$ws.Range('E17').Value = '  +7.68%  '
$ws.Range('D18').Value = '26.512.59'
$ws.Range('D19').Value = '0.9998'
$ws.Range('E19').Value = '  -0.09%  '
$ws.Range('E20').Value = '  -1.70%  '
$ws.Range('D21').Value = '1.950.92'
$ws.Range('E21').Value = '  -0.45%  '
$ws.Range('D22').Value = '4.576'
$ws.Range('E22').Value = '  +0.36%  '
$ws.Range('E23').Value = '  -1.81%  '
$ws.Range('D24').Value = '5.235'
$ws.Range('E24').Value = '  -1.78%  '
$ws.Range('D25').Value = '140.81'
$ws.Range('E25').Value = '  +3.67%  '
$ws.Range('E26').Value = '  +0.76%  '
$ws.Range('E27').Value = '  +1.12%  '
$ws.Range('B28').Value = 'BitcoinCash'
$ws.Range('C28').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D28').Value = '108.07'
$ws.Range('E28').Value = '  +1.52%  '
$ws.Range('B29').Value = 'LidoDAOToken'
$ws.Range('C29').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D29').Value = '1.767'
$ws.Range('E29').Value = '  -1.70%  '
$ws.Range('D30').Value = '4.034'
$ws.Range('E30').Value = '  +1.56%  '
$ws.Range('D31').Value = '0.08042'
$ws.Range('E31').Value = '  +2.06%  '
$ws.Range('D32').Value = '3.718'
$ws.Range('E32').Value = '  +0.34%  '
$ws.Range('D33').Value = '0.04584'
$ws.Range('E33').Value = '  +0.56%  '
$ws.Range('B34').Value = 'Frax'
$ws.Range('C34').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('D34').Value = '0.9995'
$ws.Range('E34').Value = '  -0.09%  '
$ws.Range('B35').Value = 'HuobiToken'
$ws.Range('C35').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D35').Value = '2.610'
$ws.Range('E35').Value = '  -0.23%  '
$ws.Range('B36').Value = 'ARBITRUM'
$ws.Range('C36').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D36').Value = '1.009'
$ws.Range('E36').Value = '  +1.74%  '
$ws.Range('B37').Value = 'ImmutableX'
$ws.Range('C37').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D37').Value = '0.6367'
$ws.Range('E37').Value = '  +0.36%  '
$ws.Range('B38').Value = 'TrustWalletToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D38').Value = '0.8954'
$ws.Range('E38').Value = '  -4.05%  '
$ws.Range('B39').Value = 'RenderToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D39').Value = '2.014'
$ws.Range('E39').Value = '  +1.74%  '
$ws.Range('B40').Value = 'MXToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D40').Value = '2.400'
$ws.Range('E40').Value = '  -1.27%  '
$ws.Range('B41').Value = 'PaxDollar'
$ws.Range('C41').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D41').Value = '1.004'
$ws.Range('E41').Value = '  -0.10%  '
$ws.Range('D42').Value = '103.10'
$ws.Range('E42').Value = '  -7.05%  '
$ws.Range('B43').Value = 'VeChain'
$ws.Range('C43').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D43').Value = '0.01503'
$ws.Range('E43').Value = '  -0.42%  '
$ws.Range('B44').Value = 'FraxShare'
$ws.Range('C44').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D44').Value = '5.400'
$ws.Range('E44').Value = '  -5.23%  '
$ws.Range('B45').Value = 'TheSandbox'
$ws.Range('C45').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D45').Value = '0.3890'
$ws.Range('E45').Value = '  -0.25%  '
$ws.Range('B46').Value = 'Aptos'
$ws.Range('C46').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D46').Value = '6.889'
$ws.Range('E46').Value = '  -0.23%  '
$ws.Range('B47').Value = 'Algorand'
$ws.Range('C47').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D47').Value = '0.1185'
$ws.Range('E47').Value = '  -0.59%  '
$ws.Range('B48').Value = 'Cronos'
$ws.Range('C48').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D48').Value = '0.05392'
$ws.Range('E48').Value = '  +1.08%  '
$ws.Range('B49').Value = 'Elrond'
$ws.Range('C49').Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range('D49').Value = '30.51'
$ws.Range('E49').Value = '  -0.86%  '
$ws.Range('B50').Value = 'EnergySwap'
$ws.Range('C50').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D50').Value = '7.794'
$ws.Range('E50').Value = '  -0.96%  '
$ws.Range('B51').Value = 'NEARProtocol'
$ws.Range('C51').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D51').Value = '1.262'
$ws.Range('E51').Value = '  +0.53%  '
